# Insert a new weekly price-report row for "Coliflor" (Femacal de La Calera)
# right after the existing row 314, pushing the old rows 315-347 down to
# 316-348 (dimension grows from A1:R347 to A1:R348).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 315:347 down by inserting a blank row at 315.
$ws.Rows.Item(315).Insert()

# Populate the newly inserted row 315 with the new record.
$ws.Cells.Item(315, 1).Value  = 3
$ws.Cells.Item(315, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(315, 3).Value  = "Coquimbo"
$ws.Cells.Item(315, 4).Value  = 44449
$ws.Cells.Item(315, 5).Value  = 5
$ws.Cells.Item(315, 6).Value  = 100112008
$ws.Cells.Item(315, 7).Value  = "Coliflor"
$ws.Cells.Item(315, 8).Value  = "Sin especificar"
$ws.Cells.Item(315, 9).Value  = "Primera"
$ws.Cells.Item(315, 10).Value = 1600
$ws.Cells.Item(315, 11).Value = 600
$ws.Cells.Item(315, 12).Value = 600
$ws.Cells.Item(315, 13).Value = 600
$ws.Cells.Item(315, 14).Value = "$/unidad"
$ws.Cells.Item(315, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(315, 16).Value = 600
$ws.Cells.Item(315, 17).Value = 1
$ws.Cells.Item(315, 18).Value = "Hortaliza"
